$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.962.10"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "1.635.16"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'214.42"
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -1.49%  "
$ws.Range("E9").Value = "  -2.56%  "
$ws.Range("D10").Value = "'18.53"
$ws.Range("E10").Value = "  -5.73%  "
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("D12").Value = "1.864.10"
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("D13").Value = "1.615.92"
$ws.Range("E13").Value = "  -1.93%  "
$ws.Range("D14").Value = "'4.18"
$ws.Range("E14").Value = "  -2.39%  "
$ws.Range("D15").Value = "'0.532"
$ws.Range("E15").Value = "  -2.47%  "
$ws.Range("D16").Value = "25.976.36"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("D17").Value = "0.0₃0744"
$ws.Range("E17").Value = "  -2.57%  "
$ws.Range("D18").Value = "'61.72"
$ws.Range("E18").Value = "  -2.31%  "
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "'190.74"
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("D21").Value = "'4.24"
$ws.Range("E21").Value = "  -2.36%  "
$ws.Range("D22").Value = "'9.66"
$ws.Range("E22").Value = "  -2.85%  "
$ws.Range("D23").Value = "'6.11"
$ws.Range("E23").Value = "  -2.13%  "
$ws.Range("D24").Value = "'0.134"
$ws.Range("E24").Value = "  +1.21%  "
$ws.Range("D25").Value = "'143.60"
$ws.Range("E25").Value = "  -0.50%  "
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "'6.83"
$ws.Range("E28").Value = "  -1.24%  "
$ws.Range("D29").Value = "'15.26"
$ws.Range("E29").Value = "  -1.75%  "
$ws.Range("E30").Value = "  -1.45%  "
$ws.Range("D31").Value = "'0.0483"
$ws.Range("E31").Value = "  -3.19%  "
$ws.Range("E32").Value = "  -2.80%  "
$ws.Range("E33").Value = "  -4.08%  "
$ws.Range("E34").Value = "  -2.18%  "
$ws.Range("E35").Value = "  -2.61%  "
$ws.Range("D36").Value = "1.136.32"
$ws.Range("E36").Value = "  +0.74%  "
$ws.Range("D37").Value = "'0.866"
$ws.Range("E37").Value = "  -4.34%  "
$ws.Range("E38").Value = "  -1.11%  "
$ws.Range("D39").Value = "'0.522"
$ws.Range("E39").Value = "  -3.29%  "
$ws.Range("E40").Value = "  -1.13%  "
$ws.Range("D41").Value = "'98.47"
$ws.Range("E41").Value = "  -0.77%  "
$ws.Range("D42").Value = "'0.779"
$ws.Range("E42").Value = "  -2.21%  "
$ws.Range("D43").Value = "'5.24"
$ws.Range("E43").Value = "  -4.75%  "
$ws.Range("D44").Value = "1.773.69"
$ws.Range("E44").Value = "  -0.53%  "
$ws.Range("D45").Value = "0.0₆0115"
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("D46").Value = "'55.30"
$ws.Range("E46").Value = "  -2.23%  "
$ws.Range("D47").Value = "'0.0529"
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("E48").Value = "  +2.89%  "
$ws.Range("E50").Value = "  -1.58%  "
$ws.Range("E51").Value = "  -0.01%  "
